# Generate Report for Handback
#
# This simulates a localization "handback" event: two files
# (be7a9456-dcd1-49cc-9732-4ad1a81d3fa1.md and
#  d8689d0c-8d6a-4d3d-8533-bfd0da39abbc.md) move from status
# "Ready for handoff" to "Handed back: in sync with en-US" for both the
# zh-cn and de-de locales, and their "Latest Target File" / "Latest
# Handback File" / "Latest Handback DateTime" columns get populated.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"
$hyperlinkColor = 15570276  # RGB(100,149,237) == FF6495ED, the workbook's custom hyperlink color

# ---------------------------------------------------------------------
# Sheet "Overview": update zh-cn / de-de status columns (E, F) for the
# two rows, for both be7a9456 (row 4) and d8689d0c (row 5)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(4, 5).Value = $statusHandedBack
$wsOverview.Cells.Item(4, 6).Value = $statusHandedBack
$wsOverview.Cells.Item(5, 5).Value = $statusHandedBack
$wsOverview.Cells.Item(5, 6).Value = $statusHandedBack

# ---------------------------------------------------------------------
# Sheet "zh-cn": row 4 (be7a9456) and row 5 (d8689d0c)
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 4 - be7a9456
$wsZh.Cells.Item(4, 3).Value = $statusHandedBack
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(4, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/a2fb409b44e9f091072ce1c9a71d8cdba5afdd4c/e2e/be7a9456-dcd1-49cc-9732-4ad1a81d3fa1.md", "", "", "be7a9456-dcd1-49cc-9732-4ad1a81d3fa1.md") | Out-Null
$wsZh.Cells.Item(4, 9).Font.Underline = 2
$wsZh.Cells.Item(4, 9).Font.Color = $hyperlinkColor
$wsZh.Cells.Item(4, 10).Value = "be7a9456-dcd1-49cc-9732-4ad1a81d3fa1.03c1e85547b496a5c92ff86ca1efe0b88a94cf10.zh-cn.xlf"
$wsZh.Cells.Item(4, 11).Value = "2016-08-16 04:26:15"

# Row 5 - d8689d0c
$wsZh.Cells.Item(5, 3).Value = $statusHandedBack
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(5, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/a2fb409b44e9f091072ce1c9a71d8cdba5afdd4c/e2e/d8689d0c-8d6a-4d3d-8533-bfd0da39abbc.md", "", "", "d8689d0c-8d6a-4d3d-8533-bfd0da39abbc.md") | Out-Null
$wsZh.Cells.Item(5, 9).Font.Underline = 2
$wsZh.Cells.Item(5, 9).Font.Color = $hyperlinkColor
$wsZh.Cells.Item(5, 10).Value = "d8689d0c-8d6a-4d3d-8533-bfd0da39abbc.0ad900967b0644d168275665a96eb762b92becbb.zh-cn.xlf"
$wsZh.Cells.Item(5, 11).Value = "2016-08-16 04:26:15"

# ---------------------------------------------------------------------
# Sheet "de-de": row 4 (be7a9456) and row 5 (d8689d0c)
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 4 - be7a9456
$wsDe.Cells.Item(4, 3).Value = $statusHandedBack
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(4, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/244c798cbe912f5b7b70115eeb45036211aa4eca/e2e/be7a9456-dcd1-49cc-9732-4ad1a81d3fa1.md", "", "", "be7a9456-dcd1-49cc-9732-4ad1a81d3fa1.md") | Out-Null
$wsDe.Cells.Item(4, 9).Font.Underline = 2
$wsDe.Cells.Item(4, 9).Font.Color = $hyperlinkColor
$wsDe.Cells.Item(4, 10).Value = "be7a9456-dcd1-49cc-9732-4ad1a81d3fa1.03c1e85547b496a5c92ff86ca1efe0b88a94cf10.de-de.xlf"
$wsDe.Cells.Item(4, 11).Value = "2016-08-16 04:26:22"

# Row 5 - d8689d0c
$wsDe.Cells.Item(5, 3).Value = $statusHandedBack
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(5, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/244c798cbe912f5b7b70115eeb45036211aa4eca/e2e/d8689d0c-8d6a-4d3d-8533-bfd0da39abbc.md", "", "", "d8689d0c-8d6a-4d3d-8533-bfd0da39abbc.md") | Out-Null
$wsDe.Cells.Item(5, 9).Font.Underline = 2
$wsDe.Cells.Item(5, 9).Font.Color = $hyperlinkColor
$wsDe.Cells.Item(5, 10).Value = "d8689d0c-8d6a-4d3d-8533-bfd0da39abbc.0ad900967b0644d168275665a96eb762b92becbb.de-de.xlf"
$wsDe.Cells.Item(5, 11).Value = "2016-08-16 04:26:22"
